$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7, shifting the existing rows 7-20 down to 8-21.
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with the new weekly record.
$ws.Cells.Item(7, 1).Value = 1
$ws.Cells.Item(7, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(7, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(7, 4).Value = 45002
$ws.Cells.Item(7, 5).Value = 15
$ws.Cells.Item(7, 6).Value = "Fruta"
$ws.Cells.Item(7, 7).Value = 100101
$ws.Cells.Item(7, 8).Value = "Berries"
$ws.Cells.Item(7, 9).Value = 100101007
$ws.Cells.Item(7, 10).Value = "Kiwi"
$ws.Cells.Item(7, 11).Value = "Hayward"
$ws.Cells.Item(7, 12).Value = "Segunda"
$ws.Cells.Item(7, 13).Value = 300
$ws.Cells.Item(7, 14).Value = 24000
$ws.Cells.Item(7, 15).Value = 25000
$ws.Cells.Item(7, 16).Value = 24500
$ws.Cells.Item(7, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(7, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(7, 19).Value = 1361
$ws.Cells.Item(7, 20).Value = 18
